$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, $text)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

# Row 2
$ws.Range("D2").Value = '43.260.12'
$ws.Range("E2").Value = '  -1.06%  '

# Row 3
$ws.Range("D3").Value = '2.358.14'
$ws.Range("E3").Value = '  +4.65%  '

# Row 4
$ws.Range("E4").Value = '  -0.12%  '

# Row 5
Set-TextValue $ws.Range("D5") '233.60'
$ws.Range("E5").Value = '  +0.78%  '

# Row 6
Set-TextValue $ws.Range("D6") '0.651'
$ws.Range("E6").Value = '  +0.37%  '

# Row 7
Set-TextValue $ws.Range("D7") '71.56'
$ws.Range("E7").Value = '  +12.64%  '

# Row 8
$ws.Range("E8").Value = '  -0.01%  '

# Row 9
Set-TextValue $ws.Range("D9") '0.484'
$ws.Range("E9").Value = '  +9.06%  '

# Row 10
Set-TextValue $ws.Range("D10") '0.0979'
$ws.Range("E10").Value = '  +0.93%  '

# Row 11
Set-TextValue $ws.Range("D11") '27.14'
$ws.Range("E11").Value = '  +2.38%  '

# Row 12
$ws.Range("B12").Value = 'TRON'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextValue $ws.Range("D12") '0.107'
$ws.Range("E12").Value = '  +1.96%  '

# Row 13
$ws.Range("B13").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C13").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D13").Value = '2.707.39'
$ws.Range("E13").Value = '  +4.59%  '

# Row 14
Set-TextValue $ws.Range("D14") '16.16'
$ws.Range("E14").Value = '  +3.58%  '

# Row 15
Set-TextValue $ws.Range("D15") '6.30'
$ws.Range("E15").Value = '  +3.17%  '

# Row 16
$ws.Range("E16").Value = '  +2.85%  '

# Row 17
$ws.Range("D17").Value = '2.343.78'
$ws.Range("E17").Value = '  +3.64%  '

# Row 18
$ws.Range("D18").Value = '43.301.15'

# Row 19
$ws.Range("E19").Value = '  +3.87%  '

# Row 20
$ws.Range("E20").Value = '  +3.20%  '

# Row 21
Set-TextValue $ws.Range("D21") '74.43'
$ws.Range("E21").Value = '  +1.13%  '

# Row 22
Set-TextValue $ws.Range("D22") '249.91'
$ws.Range("E22").Value = '  +0.68%  '

# Row 23
$ws.Range("B23").Value = 'WEMIXToken'
$ws.Range("C23").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue $ws.Range("D23") '3.81'
$ws.Range("E23").Value = '  +3.97%  '

# Row 24
$ws.Range("B24").Value = 'Dai'
$ws.Range("C24").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue $ws.Range("D24") '1.00'
$ws.Range("E24").Value = '  -0.06%  '

# Row 25
Set-TextValue $ws.Range("D25") '2.45'
$ws.Range("E25").Value = '  +0.17%  '

# Row 26
Set-TextValue $ws.Range("D26") '10.03'
$ws.Range("E26").Value = '  +1.13%  '

# Row 27
$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue $ws.Range("D27") '22.55'
$ws.Range("E27").Value = '  +3.81%  '

# Row 28
$ws.Range("B28").Value = 'Toncoin'
$ws.Range("C28").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue $ws.Range("D28") '2.19'
$ws.Range("E28").Value = '  -4.83%  '

# Row 29
Set-TextValue $ws.Range("D29") '172.96'
$ws.Range("E29").Value = '  -0.15%  '

# Row 30
$ws.Range("E30").Value = '  +6.08%  '

# Row 31
$ws.Range("E31").Value = '  -4.15%  '

# Row 32
Set-TextValue $ws.Range("D32") '0.128'
$ws.Range("E32").Value = '  +0.98%  '

# Row 33
$ws.Range("E33").Value = '  +1.25%  '

# Row 34
Set-TextValue $ws.Range("D34") '0.0690'
$ws.Range("E34").Value = '  +1.39%  '

# Row 35
Set-TextValue $ws.Range("D35") '5.07'
$ws.Range("E35").Value = '  +2.90%  '

# Row 36
$ws.Range("B36").Value = 'RenderToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws.Range("D36") '3.70'
$ws.Range("E36").Value = '  +1.65%  '

# Row 37
$ws.Range("B37").Value = 'THORChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
Set-TextValue $ws.Range("D37") '6.55'
$ws.Range("E37").Value = '  +3.01%  '

# Row 38
$ws.Range("E38").Value = '  +6.06%  '

# Row 39
$ws.Range("E39").Value = '  +0.02%  '

# Row 40
$ws.Range("E40").Value = '  -0.01%  '

# Row 41
Set-TextValue $ws.Range("D41") '8.91'
$ws.Range("E41").Value = '  +3.73%  '

# Row 42
Set-TextValue $ws.Range("D42") '18.64'
$ws.Range("E42").Value = '  +8.67%  '

# Row 43
$ws.Range("E43").Value = '  +8.08%  '

# Row 44
Set-TextValue $ws.Range("D44") '99.17'
$ws.Range("E44").Value = '  +0.99%  '

# Row 45
Set-TextValue $ws.Range("D45") '4.48'
$ws.Range("E45").Value = '  -3.26%  '

# Row 46
$ws.Range("E46").Value = '  +2.04%  '

# Row 47
Set-TextValue $ws.Range("D47") '0.0953'
$ws.Range("E47").Value = '  +1.25%  '

# Row 48
$ws.Range("D48").Value = '1.441.12'
$ws.Range("E48").Value = '  -0.74%  '

# Row 49
$ws.Range("D49").Value = '2.583.40'
$ws.Range("E49").Value = '  +4.86%  '

# Row 50
$ws.Range("B50").Value = 'TerraClassic'
$ws.Range("C50").Value = 'https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc'
Set-TextValue $ws.Range("D50") '0.000203'
$ws.Range("E50").Value = '  -2.36%  '

# Row 51
$ws.Range("B51").Value = 'HuobiToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextValue $ws.Range("D51") '2.77'
$ws.Range("E51").Value = '  +0.33%  '
